# Compatibilização com novo formato da base da ANEEL
# Adds six new sig_agente / nome_4md mapping rows to the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("CELESC", "CELESC"),
    @("EQUATORIAL MA", "EQUATORIAL MA"),
    @("EQUATORIAL AL", "EQUATORIAL AL"),
    @("BOA VISTA", "RORAIMA"),
    @("CERAL ARARUAMA", "OUTRA"),
    @("CERSAD DISTRIBUI", "OUTRA")
)

$startRow = 134
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Reflect the cursor/selection position shown in the edited file (cosmetic).
$ws.Range("B137").Select() | Out-Null
